# Apply crypto list updates (prices and 1h volume %) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.267.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.713.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.89%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "658.59"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.435"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.21%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("E10").Value = "  -2.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.710.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000319"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +17.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.01%  "

$ws.Range("E14").Value = "  +0.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.407.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.762.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.702.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.510"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "526.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.07%  "

$ws.Range("E24").Value = "  -0.92%  "

$ws.Range("E25").Value = "  +9.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "106.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.30%  "

$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.43%  "

$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.190"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +12.78%  "

$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.912.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.56%  "

$ws.Range("E32").Value = "  -0.84%  "

$ws.Range("E33").Value = "  -0.18%  "

$ws.Range("E34").Value = "  +3.00%  "

$ws.Range("E35").Value = "  -4.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "641.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.592"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.05%  "

$ws.Range("E40").Value = "  -1.47%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.167"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.485"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.64%  "

$ws.Range("E48").Value = "  -1.53%  "

$ws.Range("E49").Value = "  +2.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.23%  "
